$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain plain text so numeric-looking price strings
# (e.g. "1.00", "174.65") are not coerced into floating point numbers,
# matching the workbook's original inline-string "Price" column.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.968.61"
$ws.Range("E2").Value = "  -2.87%  "

$ws.Range("D3").Value = "3.488.69"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "583.69"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("D6").Value = "174.65"
$ws.Range("E6").Value = "  -3.64%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  -3.17%  "

$ws.Range("D9").Value = "3.486.85"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("E10").Value = "  -6.40%  "

$ws.Range("E11").Value = "  -2.10%  "

$ws.Range("D13").Value = "4.086.73"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("D15").Value = "30.24"
$ws.Range("E15").Value = "  -5.78%  "

$ws.Range("D16").Value = "66.056.75"
$ws.Range("E16").Value = "  -2.82%  "

$ws.Range("E17").Value = "  -2.88%  "

$ws.Range("D18").Value = "3.485.08"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("E19").Value = "  -4.29%  "

$ws.Range("D20").Value = "13.96"
$ws.Range("E20").Value = "  -0.90%  "

$ws.Range("D21").Value = "367.14"
$ws.Range("E21").Value = "  -6.67%  "

$ws.Range("E22").Value = "  -1.94%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "72.55"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  +4.52%  "

$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  -7.23%  "

$ws.Range("E28").Value = "  +0.99%  "

$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.26%  "

$ws.Range("D30").Value = "24.20"
$ws.Range("E30").Value = "  +2.48%  "

$ws.Range("D31").Value = "5.79"
$ws.Range("E31").Value = "  -5.33%  "

$ws.Range("E32").Value = "  -3.46%  "

$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").Value = "7.16"
$ws.Range("E34").Value = "  -2.87%  "

$ws.Range("E35").Value = "  -7.42%  "

$ws.Range("E36").Value = "  -1.58%  "

$ws.Range("D37").Value = "160.11"
$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("D38").Value = "29.46"
$ws.Range("E38").Value = "  +13.00%  "

$ws.Range("D39").Value = "0.890"
$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("D40").Value = "2.836.63"
$ws.Range("E40").Value = "  +3.67%  "

$ws.Range("E41").Value = "  -5.45%  "

$ws.Range("D42").Value = "2.60"
$ws.Range("E42").Value = "  -7.19%  "

$ws.Range("D43").Value = "4.44"
$ws.Range("E43").Value = "  -4.55%  "

$ws.Range("D44").Value = "6.43"
$ws.Range("E44").Value = "  -4.49%  "

$ws.Range("D45").Value = "0.0684"
$ws.Range("E45").Value = "  -4.66%  "

$ws.Range("D46").Value = "39.93"
$ws.Range("E46").Value = "  -3.50%  "

$ws.Range("D47").Value = "24.29"
$ws.Range("E47").Value = "  -7.08%  "

$ws.Range("E48").Value = "  -3.48%  "

$ws.Range("D49").Value = "309.49"
$ws.Range("E49").Value = "  -5.85%  "

$ws.Range("D50").Value = "0.824"
$ws.Range("E50").Value = "  -2.47%  "

$ws.Range("D51").Value = "6.21"
$ws.Range("E51").Value = "  -2.47%  "
